# "added data driven concept using Excel"
# Remove the "email" column and refresh the registration sample rows to a
# new data set (testuser5..testuser8). The phone numbers move from column D
# into column C (now text, entered with a leading apostrophe so Excel keeps
# the quote-prefix / text formatting instead of re-parsing them as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the whole "email" column (C); phone/password/occupation/gender shift left.
$ws.Range("C1:C5").EntireColumn.Delete()

# --- Row 2 (was testuser1) ---
$ws.Range("A2").Value = "testuser5"
$ws.Range("B2").Value = "lastname5"
$ws.Range("C2").Value = "'1234567898"
$ws.Range("D2").Value = "Selenium@123"
$ws.Range("E2").Value = "Student"
$ws.Range("F2").Value = "male"

# --- Row 3 (was testuser2) ---
$ws.Range("A3").Value = "testuser6"
$ws.Range("B3").Value = "lastname6"
$ws.Range("C3").Value = "'2234567898"
$ws.Range("D3").Value = "Selenium@124"
$ws.Range("E3").Value = "Doctor"
$ws.Range("F3").Value = "female"

# --- Row 4 (was testuser3) ---
$ws.Range("A4").Value = "testuser7"
$ws.Range("B4").Value = "lastname7"
$ws.Range("C4").Value = "'3123456789"
$ws.Range("D4").Value = "Selenium@125"
$ws.Range("E4").Value = "Engineer"
$ws.Range("F4").Value = "male"

# --- Row 5 (was testuser4) ---
$ws.Range("A5").Value = "testuser8"
$ws.Range("B5").Value = "lastname8"
$ws.Range("C5").Value = "'4234567898"
$ws.Range("D5").Value = "Selenium@126"
$ws.Range("E5").Value = "Scientist"
$ws.Range("F5").Value = "female"
